# "Début de transformation du model, maintenant en model relationnel de BD"
# Mark the next batch of checklist rows (13,14,17,18,19,20,21 in column C)
# as done ("X"), same convention already used for rows 4-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

foreach ($r in 13, 14, 17, 18, 19, 20, 21) {
    $ws.Cells.Item($r, 3).Value = "X"
}

# Scroll the frozen view down so row 6 becomes the first visible row under
# the frozen header row, and move the live selection to I23.
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("I23").Select()
